$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7409614
$ws.Range("I51").Value = 2599.75
$ws.Range("J51").Value = 10103074
$ws.Range("K51").Value = 2599.75
$ws.Range("L51").Value = 10103074
$ws.Range("M51").Value = -2115.75
$ws.Range("N51").Value = -10104042
$ws.Range("H113").Value = 2275.2
$ws.Range("I113").Value = 1630.5
$ws.Range("K113").Value = 1630.5
$ws.Range("M113").Value = 1623.5
$ws.Range("H132").Value = 5143.1904
$ws.Range("J132").Value = 3361.6
$ws.Range("L132").Value = 10084.8
$ws.Range("N132").Value = -15144.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1920.1578
$ws.Range("I110").Value = 1544
$ws.Range("J110").Value = 2735.1667
$ws.Range("K110").Value = 1544
$ws.Range("L110").Value = 2735.1667
$ws.Range("M110").Value = 501
$ws.Range("N110").Value = -6825.1667
$ws.Range("H122").Value = 3764.6667
$ws.Range("I122").Value = 850
$ws.Range("J122").Value = 5222
$ws.Range("K122").Value = 2550
$ws.Range("L122").Value = 15666
$ws.Range("M122").Value = -100
$ws.Range("N122").Value = -20566
$ws.Range("H132").Value = 2987.4285
$ws.Range("I132").Value = 1943.5
$ws.Range("J132").Value = 3936.4546
$ws.Range("K132").Value = 5830.5
$ws.Range("L132").Value = 11809.3638
$ws.Range("M132").Value = -3300.5
$ws.Range("N132").Value = -16869.3638
$ws.Range("H134").Value = 25311.6
$ws.Range("J134").Value = 25311.6
$ws.Range("L134").Value = 25311.6
$ws.Range("N134").Value = -35451.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1500.7333
$ws.Range("I107").Value = 1003.875
$ws.Range("J107").Value = 2068.5715
$ws.Range("K107").Value = 1003.875
$ws.Range("L107").Value = 2068.5715
$ws.Range("M107").Value = 916.125
$ws.Range("N107").Value = -5908.5715
$ws.Range("H134").Value = 1513.3636
$ws.Range("I134").Value = 1195.4
$ws.Range("J134").Value = 2002.5385
$ws.Range("K134").Value = 3586.2
$ws.Range("L134").Value = 6007.6155
$ws.Range("M134").Value = -1051.2
$ws.Range("N134").Value = -11077.6155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1398.5769
$ws.Range("I58").Value = 823.8570999999999
$ws.Range("J58").Value = 1610.3158
$ws.Range("K58").Value = 823.8570999999999
$ws.Range("L58").Value = 1610.3158
$ws.Range("M58").Value = -620.8570999999999
$ws.Range("N58").Value = -2016.3158
$ws.Range("H107").Value = 21740994
$ws.Range("I107").Value = 71430664
$ws.Range("J107").Value = 1762.5
$ws.Range("K107").Value = 71430664
$ws.Range("L107").Value = 1762.5
$ws.Range("M107").Value = -71428744
$ws.Range("N107").Value = -5602.5
$ws.Range("H122").Value = 5556479
$ws.Range("I122").Value = 6667433.5
$ws.Range("J122").Value = 1707
$ws.Range("K122").Value = 20002300.5
$ws.Range("L122").Value = 5121
$ws.Range("M122").Value = -19999850.5
$ws.Range("N122").Value = -10021
$ws.Range("H132").Value = 2053.4348
$ws.Range("I132").Value = 1635.8334
$ws.Range("J132").Value = 2509
$ws.Range("K132").Value = 4907.5002
$ws.Range("L132").Value = 7527
$ws.Range("M132").Value = -2377.5002
$ws.Range("N132").Value = -12587
$ws.Range("H134").Value = 4867.1
$ws.Range("I134").Value = 4896.375
$ws.Range("J134").Value = 4750
$ws.Range("K134").Value = 14689.125
$ws.Range("L134").Value = 14250
$ws.Range("M134").Value = -12154.125
$ws.Range("N134").Value = -19320
$ws.Range("H136").Value = 1398.5769
$ws.Range("I136").Value = 823.8570999999999
$ws.Range("J136").Value = 1610.3158
$ws.Range("K136").Value = 2471.5713
$ws.Range("L136").Value = 4830.9474
$ws.Range("M136").Value = 78.42870000000039
$ws.Range("N136").Value = -9930.947400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 152384.16
$ws.Range("I2").Value = 282907
$ws.Range("J2").Value = 107.5
$ws.Range("K2").Value = 1697442
$ws.Range("L2").Value = 645
$ws.Range("M2").Value = -1697329
$ws.Range("N2").Value = -871
$ws.Range("H68").Value = 551.93335
$ws.Range("I68").Value = 375.9375
$ws.Range("J68").Value = 753.0714
$ws.Range("K68").Value = 1127.8125
$ws.Range("L68").Value = 2259.2142
$ws.Range("M68").Value = -316.8125
$ws.Range("N68").Value = -3881.2142
$ws.Range("H71").Value = 551.93335
$ws.Range("I71").Value = 375.9375
$ws.Range("J71").Value = 753.0714
$ws.Range("K71").Value = 3383.4375
$ws.Range("L71").Value = 6777.6426
$ws.Range("M71").Value = 672.5625
$ws.Range("N71").Value = -14889.6426
$ws.Range("H86").Value = 194.33333
$ws.Range("I86").Value = 193.2
$ws.Range("J86").Value = 200
$ws.Range("K86").Value = 579.5999999999999
$ws.Range("L86").Value = 600
$ws.Range("M86").Value = 606.4000000000001
$ws.Range("N86").Value = -2972
$ws.Range("H89").Value = 194.33333
$ws.Range("I89").Value = 193.2
$ws.Range("J89").Value = 200
$ws.Range("K89").Value = 1738.8
$ws.Range("L89").Value = 1800
$ws.Range("M89").Value = 4189.2
$ws.Range("N89").Value = -13656
$ws.Range("H133").Value = 1069.75
$ws.Range("I133").Value = 747.4
$ws.Range("K133").Value = 2242.2
$ws.Range("M133").Value = 2817.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 57147116
$ws.Range("I70").Value = 88892830
$ws.Range("K70").Value = 88892830
$ws.Range("M70").Value = -88892560
$ws.Range("H73").Value = 57147116
$ws.Range("I73").Value = 88892830
$ws.Range("K73").Value = 88892830
$ws.Range("M73").Value = -88891894
$ws.Range("H80").Value = 2564.8965
$ws.Range("I80").Value = 2567
$ws.Range("J80").Value = 2564.348
$ws.Range("K80").Value = 2567
$ws.Range("L80").Value = 2564.348
$ws.Range("M80").Value = -1569
$ws.Range("N80").Value = -4560.348
$ws.Range("H83").Value = 2564.8965
$ws.Range("I83").Value = 2567
$ws.Range("J83").Value = 2564.348
$ws.Range("K83").Value = 12835
$ws.Range("L83").Value = 12821.74
$ws.Range("M83").Value = -7843
$ws.Range("N83").Value = -22805.74
$ws.Range("H88").Value = 39750
$ws.Range("J88").Value = 39750
$ws.Range("L88").Value = 39750
$ws.Range("N88").Value = -40652
$ws.Range("H91").Value = 39750
$ws.Range("J91").Value = 39750
$ws.Range("L91").Value = 39750
$ws.Range("N91").Value = -42870
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 3128.5
$ws.Range("I122").Value = 3871.3333
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 11613.9999
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -9163.999899999999
$ws.Range("N122").Value = -7600
$ws.Range("H132").Value = 3999.25
$ws.Range("I132").Value = 3666
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 10998
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -8468
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1429.6666
$ws.Range("I40").Value = 1480.6364
$ws.Range("J40").Value = 1386.5385
$ws.Range("K40").Value = 1480.6364
$ws.Range("L40").Value = 1386.5385
$ws.Range("M40").Value = -1344.6364
$ws.Range("N40").Value = -1658.5385
$ws.Range("H122").Value = 2369.4783
$ws.Range("I122").Value = 2295.1765
$ws.Range("J122").Value = 2580
$ws.Range("K122").Value = 6885.529500000001
$ws.Range("L122").Value = 7740
$ws.Range("M122").Value = -4435.529500000001
$ws.Range("N122").Value = -12640
$ws.Range("H132").Value = 14294079
$ws.Range("I132").Value = 22738836
$ws.Range("J132").Value = 2952
$ws.Range("K132").Value = 68216508
$ws.Range("L132").Value = 8856
$ws.Range("M132").Value = -68213978
$ws.Range("N132").Value = -13916

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 85000
$ws.Range("J108").Value = 85000
$ws.Range("L108").Value = 85000
$ws.Range("N108").Value = -92680
$ws.Range("H122").Value = 1585.5714
$ws.Range("I122").Value = 1033.3334
$ws.Range("J122").Value = 1999.75
$ws.Range("K122").Value = 3100.0002
$ws.Range("L122").Value = 5999.25
$ws.Range("M122").Value = -650.0001999999999
$ws.Range("N122").Value = -10899.25
